$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 8205.773999999999
$ws.Range("I70").Value = 9347.134
$ws.Range("J70").Value = 7135.75
$ws.Range("K70").Value = 28041.402
$ws.Range("L70").Value = 21407.25
$ws.Range("M70").Value = -27771.402
$ws.Range("N70").Value = -21947.25

$ws.Range("H73").Value = 8205.773999999999
$ws.Range("I73").Value = 9347.134
$ws.Range("J73").Value = 7135.75
$ws.Range("K73").Value = 28041.402
$ws.Range("L73").Value = 21407.25
$ws.Range("M73").Value = -27105.402
$ws.Range("N73").Value = -23279.25

$ws.Range("H134").Value = 112947.62
$ws.Range("J134").Value = 107401.75
$ws.Range("L134").Value = 107401.75
$ws.Range("N134").Value = -117541.75

$ws.Range("H135").Value = 1725.9333
$ws.Range("I135").Value = 1492.1072
$ws.Range("K135").Value = 13428.9648
$ws.Range("M135").Value = -10893.9648

$ws.Range("H138").Value = 4562.8
$ws.Range("J138").Value = 4244.8936
$ws.Range("L138").Value = 12734.6808
$ws.Range("N138").Value = -23014.6808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1774.3125
$ws.Range("I2").Value = 1492.9333
$ws.Range("J2").Value = 5995
$ws.Range("K2").Value = 1492.9333
$ws.Range("L2").Value = 5995
$ws.Range("M2").Value = -1379.9333
$ws.Range("N2").Value = -6221

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H32").Value = 10028.2
$ws.Range("I32").Value = 8087.9443
$ws.Range("K32").Value = 8087.9443
$ws.Range("M32").Value = -7800.9443

$ws.Range("H63").Value = 7375
$ws.Range("J63").Value = 9750
$ws.Range("L63").Value = 9750
$ws.Range("N63").Value = -11122

$ws.Range("H66").Value = 7375
$ws.Range("J66").Value = 9750
$ws.Range("L66").Value = 48750
$ws.Range("N66").Value = -55614

$ws.Range("H102").Value = 1247.8
$ws.Range("J102").Value = 1140
$ws.Range("L102").Value = 1140
$ws.Range("N102").Value = -4384

$ws.Range("H109").Value = 100001
$ws.Range("J109").Value = 100001
$ws.Range("L109").Value = 100001
$ws.Range("N109").Value = -102775

$ws.Range("H116").Value = 1774.3125
$ws.Range("I116").Value = 1492.9333
$ws.Range("J116").Value = 5995
$ws.Range("K116").Value = 1492.9333
$ws.Range("L116").Value = 5995
$ws.Range("M116").Value = 801.0667000000001
$ws.Range("N116").Value = -10583

$ws.Range("H132").Value = 9697.267
$ws.Range("I132").Value = 11770.591
$ws.Range("K132").Value = 35311.773
$ws.Range("M132").Value = -32781.773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1774.3125
$ws.Range("I3").Value = 1492.9333
$ws.Range("J3").Value = 5995
$ws.Range("K3").Value = 1492.9333
$ws.Range("L3").Value = 5995
$ws.Range("M3").Value = -1378.9333
$ws.Range("N3").Value = -6223

$ws.Range("H80").Value = 3204.077
$ws.Range("I80").Value = 734.3333
$ws.Range("J80").Value = 3945
$ws.Range("K80").Value = 734.3333
$ws.Range("L80").Value = 3945
$ws.Range("M80").Value = 263.6667
$ws.Range("N80").Value = -5941

$ws.Range("H82").Value = 21360.818
$ws.Range("J82").Value = 24996.25
$ws.Range("L82").Value = 24996.25
$ws.Range("N82").Value = -25762.25

$ws.Range("H83").Value = 3204.077
$ws.Range("I83").Value = 734.3333
$ws.Range("J83").Value = 3945
$ws.Range("K83").Value = 3671.6665
$ws.Range("L83").Value = 19725
$ws.Range("M83").Value = 1320.3335
$ws.Range("N83").Value = -29709

$ws.Range("H85").Value = 21360.818
$ws.Range("J85").Value = 24996.25
$ws.Range("L85").Value = 24996.25
$ws.Range("N85").Value = -27648.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 57494
$ws.Range("I52").Value = 54989
$ws.Range("J52").Value = 59999
$ws.Range("K52").Value = 54989
$ws.Range("L52").Value = 59999
$ws.Range("M52").Value = -54695
$ws.Range("N52").Value = -60587

$ws.Range("H58").Value = 1811.55
$ws.Range("I58").Value = 1761.6
$ws.Range("K58").Value = 1761.6
$ws.Range("M58").Value = -1558.6

$ws.Range("H134").Value = 4342.0835
$ws.Range("I134").Value = 3884.7097
$ws.Range("J134").Value = 7177.8
$ws.Range("K134").Value = 11654.1291
$ws.Range("L134").Value = 21533.4
$ws.Range("M134").Value = -9119.1291
$ws.Range("N134").Value = -26603.4

$ws.Range("H135").Value = 84840
$ws.Range("J135").Value = 84840
$ws.Range("L135").Value = 84840
$ws.Range("N135").Value = -94980

$ws.Range("H136").Value = 1811.55
$ws.Range("I136").Value = 1761.6
$ws.Range("K136").Value = 5284.799999999999
$ws.Range("M136").Value = -2734.799999999999

$ws.Range("H141").Value = 558105.5
$ws.Range("J141").Value = 597115.25
$ws.Range("L141").Value = 597115.25
$ws.Range("N141").Value = -607475.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1046.875
$ws.Range("J2").Value = 2566.8333
$ws.Range("L2").Value = 15400.9998
$ws.Range("N2").Value = -15626.9998

$ws.Range("H39").Value = 206334.6
$ws.Range("J39").Value = 7918.5
$ws.Range("L39").Value = 23755.5
$ws.Range("N39").Value = -24343.5

$ws.Range("H122").Value = 1093.65
$ws.Range("J122").Value = 1295.6666
$ws.Range("L122").Value = 11660.9994
$ws.Range("N122").Value = -16560.9994

$ws.Range("H132").Value = 2253.818
$ws.Range("J132").Value = 2313.4
$ws.Range("L132").Value = 20820.6
$ws.Range("N132").Value = -25880.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8007.231
$ws.Range("I70").Value = 9513.714
$ws.Range("K70").Value = 9513.714
$ws.Range("M70").Value = -9243.714

$ws.Range("H73").Value = 8007.231
$ws.Range("I73").Value = 9513.714
$ws.Range("K73").Value = 9513.714
$ws.Range("M73").Value = -8577.714

$ws.Range("H80").Value = 3933.3845
$ws.Range("I80").Value = 2869
$ws.Range("K80").Value = 2869
$ws.Range("M80").Value = -1871

$ws.Range("H83").Value = 3933.3845
$ws.Range("I83").Value = 2869
$ws.Range("K83").Value = 14345
$ws.Range("M83").Value = -9353

$ws.Range("H107").Value = 35715170
$ws.Range("J107").Value = 125002650
$ws.Range("L107").Value = 125002650
$ws.Range("N107").Value = -125006490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29263.846
$ws.Range("I7").Value = 42179
$ws.Range("K7").Value = 42179
$ws.Range("M7").Value = -42067

$ws.Range("H22").Value = 2161.2856
$ws.Range("I22").Value = 2381
$ws.Range("K22").Value = 2381
$ws.Range("M22").Value = -2086

$ws.Range("H27").Value = 2161.2856
$ws.Range("I27").Value = 2381
$ws.Range("K27").Value = 2381
$ws.Range("M27").Value = -2274

$ws.Range("H40").Value = 13950
$ws.Range("I40").Value = 16183.75
$ws.Range("K40").Value = 16183.75
$ws.Range("M40").Value = -16047.75

$ws.Range("H126").Value = 29263.846
$ws.Range("I126").Value = 42179
$ws.Range("K126").Value = 126537
$ws.Range("M126").Value = -124067

$ws.Range("H136").Value = 18245.277
$ws.Range("J136").Value = 7416.0835
$ws.Range("L136").Value = 22248.2505
$ws.Range("N136").Value = -27348.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2770.647
$ws.Range("I122").Value = 2889.3333
$ws.Range("J122").Value = 1880.5
$ws.Range("K122").Value = 8667.999899999999
$ws.Range("L122").Value = 5641.5
$ws.Range("M122").Value = -6217.999899999999
$ws.Range("N122").Value = -10541.5
